$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear row 4 invite data (B4, D4, F4) - keep A4 as is
$ws.Range("B4").Value = $null
$ws.Range("D4").Value = $null
$ws.Range("F4").Value = $null

# Row 5
$ws.Range("B5").Value = "Mathieu"
$ws.Range("D5").Value = "Hugo D"
$ws.Range("F5").Value = "07:37"

# Row 6
$ws.Range("B6").Value = "Mathieu"
$ws.Range("D6").Value = "Sylvie P"
$ws.Range("F6").Value = "07:32"

# Row 7
$ws.Range("B7").Value = "Steve"
$ws.Range("D7").Value = "Baptiste"
$ws.Range("F7").Value = "07:28"

# Row 8
$ws.Range("B8").Value = "Hugo W"
$ws.Range("D8").Value = "Hugo D"
$ws.Range("F8").Value = "07:28"

# Row 9
$ws.Range("B9").Value = "Hugo W"
$ws.Range("D9").Value = "Hugo D"
$ws.Range("F9").Value = $null

# Row 10
$ws.Range("B10").Value = "Baptiste"
$ws.Range("D10").Value = "Hugo D"
$ws.Range("F10").Value = $null
